# Auto-generated Excel COM-interop edit script
# Refreshes market-data derived columns (H..N: current average
# prices, leve NQ/HQ prices and profit figures) across the
# crafting-job sheets to match the latest scheduled market snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 8: On the Drip
$ws.Range("H8").Value = 199.2
$ws.Range("I8").Value = 224
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 672
$ws.Range("L8").Value = 300
$ws.Range("M8").Value = -533
$ws.Range("N8").Value = -578

# ALC row 17: One for the Road
$ws.Range("H17").Value = 3410.0217
$ws.Range("J17").Value = 3410.0217
$ws.Range("L17").Value = 10230.0651
$ws.Range("N17").Value = -10566.0651

# ALC row 108: Keeping Magic Alive
$ws.Range("H108").Value = 26860.334
$ws.Range("J108").Value = 26860.334
$ws.Range("L108").Value = 26860.334
$ws.Range("N108").Value = -34540.334

# ALC row 117: A Greater Grimoire
$ws.Range("H117").Value = 48538
$ws.Range("J117").Value = 48538
$ws.Range("L117").Value = 48538
$ws.Range("N117").Value = -57716

# ALC row 120: Supreme Official Strategy Guide
$ws.Range("H120").Value = 49657
$ws.Range("J120").Value = 49657
$ws.Range("L120").Value = 49657
$ws.Range("N120").Value = -59333

# ALC row 126: Rebuilding to Code
$ws.Range("H126").Value = 47006
$ws.Range("J126").Value = 47006
$ws.Range("L126").Value = 47006
$ws.Range("N126").Value = -56886

# ALC row 129: Practical Command
$ws.Range("H129").Value = 1329.4546
$ws.Range("I129").Value = 1531.6666
$ws.Range("J129").Value = 1253.625
$ws.Range("K129").Value = 4594.9998
$ws.Range("L129").Value = 3760.875
$ws.Range("M129").Value = 405.0002000000004
$ws.Range("N129").Value = -13760.875

# ALC row 138: All-night Crafting
$ws.Range("H138").Value = 1498.92
$ws.Range("I138").Value = 832.1429000000001
$ws.Range("J138").Value = 1857.9539
$ws.Range("K138").Value = 2496.4287
$ws.Range("L138").Value = 5573.861699999999
$ws.Range("M138").Value = 2643.5713
$ws.Range("N138").Value = -15853.8617

# ALC row 141: Remedy for Reason
$ws.Range("H141").Value = 6223.769
$ws.Range("I141").Value = 3827.8572
$ws.Range("J141").Value = 9019
$ws.Range("K141").Value = 11483.5716
$ws.Range("L141").Value = 27057
$ws.Range("M141").Value = -6303.571599999999
$ws.Range("N141").Value = -37417

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1315.0714
$ws.Range("I45").Value = 1101.5714
$ws.Range("J45").Value = 1528.5714
$ws.Range("K45").Value = 1101.5714
$ws.Range("L45").Value = 1528.5714
$ws.Range("M45").Value = -724.5714
$ws.Range("N45").Value = -2282.5714

# ARM row 109: A Head of Demand
$ws.Range("H109").Value = 41655.75
$ws.Range("J109").Value = 41655.75
$ws.Range("L109").Value = 41655.75
$ws.Range("N109").Value = -44429.75

# ARM row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1616.6
$ws.Range("I110").Value = 1696.2222
$ws.Range("K110").Value = 1696.2222
$ws.Range("M110").Value = 348.7778000000001

# ARM row 117: Signed, Shield, Delivered
$ws.Range("H117").Value = 42847.8
$ws.Range("J117").Value = 42847.8
$ws.Range("L117").Value = 42847.8
$ws.Range("N117").Value = -52025.8

# ARM row 118: A Budding Business
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314

# ARM row 138: Don't Ask about the Rivets
$ws.Range("H138").Value = 47809.668
$ws.Range("J138").Value = 47809.668
$ws.Range("L138").Value = 47809.668
$ws.Range("N138").Value = -58089.668

$ws = $wb.Worksheets.Item("BSM")
# BSM row 119: Bae Blade
$ws.Range("H119").Value = 47753
$ws.Range("J119").Value = 47753
$ws.Range("L119").Value = 47753
$ws.Range("N119").Value = -57429

# BSM row 120: Under the Fool Moon
$ws.Range("H120").Value = 48761
$ws.Range("J120").Value = 48761
$ws.Range("L120").Value = 48761
$ws.Range("N120").Value = -58437

# BSM row 126: Records of the Republic
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656

# BSM row 132: Always Be Prepaired
$ws.Range("H132").Value = 40437.5
$ws.Range("J132").Value = 40437.5
$ws.Range("L132").Value = 40437.5
$ws.Range("N132").Value = -50557.5

# BSM row 137: Dagger Swagger
$ws.Range("H137").Value = 28189.834
$ws.Range("J137").Value = 28189.834
$ws.Range("L137").Value = 28189.834
$ws.Range("N137").Value = -38389.834

$ws = $wb.Worksheets.Item("CRP")
# CRP row 20: Re-crating the Scene
$ws.Range("H20").Value = 43146.332
$ws.Range("J20").Value = 43146.332
$ws.Range("L20").Value = 43146.332
$ws.Range("N20").Value = -43618.332

# CRP row 30: Polearms Aplenty
$ws.Range("H30").Value = 43146.332
$ws.Range("J30").Value = 43146.332
$ws.Range("L30").Value = 43146.332
$ws.Range("N30").Value = -43328.332

# CRP row 99: O Pine
$ws.Range("H99").Value = 1637.8889
$ws.Range("J99").Value = 1820
$ws.Range("L99").Value = 1820
$ws.Range("N99").Value = -4816

# CRP row 115: Horde of the Rings
$ws.Range("H115").Value = 34254.332
$ws.Range("J115").Value = 34254.332
$ws.Range("L115").Value = 34254.332
$ws.Range("N115").Value = -36604.332

# CRP row 116: The Right Tool for the Job
$ws.Range("H116").Value = 49368.5
$ws.Range("J116").Value = 49368.5
$ws.Range("L116").Value = 49368.5
$ws.Range("N116").Value = -58546.5

# CRP row 126: A Better Conductor
$ws.Range("H126").Value = 1637.8889
$ws.Range("J126").Value = 1820
$ws.Range("L126").Value = 5460
$ws.Range("N126").Value = -10400

# CRP row 128: An A-prop-riate Request
$ws.Range("H128").Value = 43146.332
$ws.Range("J128").Value = 43146.332
$ws.Range("L128").Value = 43146.332
$ws.Range("N128").Value = -53106.332

# CRP row 134: Wood You Be Quiet
$ws.Range("H134").Value = 585437.2
$ws.Range("I134").Value = 1068.7059
$ws.Range("J134").Value = 2004617.9
$ws.Range("K134").Value = 3206.1177
$ws.Range("L134").Value = 6013853.699999999
$ws.Range("M134").Value = -671.1176999999998
$ws.Range("N134").Value = -6018923.699999999

# CRP row 138: Bow Out
$ws.Range("H138").Value = 46430.8
$ws.Range("J138").Value = 46430.8
$ws.Range("L138").Value = 46430.8
$ws.Range("N138").Value = -56710.8

$ws = $wb.Worksheets.Item("GSM")
# GSM row 107: Whetstones for the Workers
$ws.Range("H107").Value = 7879.8
$ws.Range("I107").Value = 999.5
$ws.Range("J107").Value = 12466.667
$ws.Range("K107").Value = 999.5
$ws.Range("L107").Value = 12466.667
$ws.Range("M107").Value = 920.5
$ws.Range("N107").Value = -16306.667

# GSM row 110: Slimming Down
$ws.Range("H110").Value = 27699.4
$ws.Range("J110").Value = 27699.4
$ws.Range("L110").Value = 27699.4
$ws.Range("N110").Value = -35879.4

# GSM row 130: Planisphere to Paper
$ws.Range("H130").Value = 46416.777
$ws.Range("J130").Value = 46416.777
$ws.Range("L130").Value = 46416.777
$ws.Range("N130").Value = -56456.777

# GSM row 134: Guaranteed Gem
$ws.Range("H134").Value = 24380
$ws.Range("J134").Value = 24380
$ws.Range("L134").Value = 73140
$ws.Range("N134").Value = -78210

# GSM row 135: Fan of the Foreign
$ws.Range("H135").Value = 31778.273
$ws.Range("J135").Value = 31778.273
$ws.Range("L135").Value = 31778.273
$ws.Range("N135").Value = -41918.273

# GSM row 136: Shiny and Good
$ws.Range("H136").Value = 20569.54
$ws.Range("J136").Value = 20569.54
$ws.Range("L136").Value = 61708.62
$ws.Range("N136").Value = -66808.62

# GSM row 138: Orders Anonymous
$ws.Range("H138").Value = 45650
$ws.Range("J138").Value = 45650
$ws.Range("L138").Value = 45650
$ws.Range("N138").Value = -55930

$ws = $wb.Worksheets.Item("LTW")
# LTW row 61: Spelling Me Softly
$ws.Range("H61").Value = 4723.1665
$ws.Range("I61").Value = 4667.8
$ws.Range("K61").Value = 4667.8
$ws.Range("M61").Value = -4465.8

# LTW row 113: Peace in Rest
$ws.Range("H113").Value = 4723.1665
$ws.Range("I113").Value = 4667.8
$ws.Range("K113").Value = 4667.8
$ws.Range("M113").Value = -2497.8

# LTW row 130: Generous Soles
$ws.Range("H130").Value = 39996
$ws.Range("J130").Value = 39996
$ws.Range("L130").Value = 39996
$ws.Range("N130").Value = -50036

# LTW row 132: Tenets of Tanning
$ws.Range("H132").Value = 3964.08
$ws.Range("I132").Value = 2237.5454
$ws.Range("K132").Value = 6712.6362
$ws.Range("M132").Value = -4182.6362

# LTW row 136: Respect for Br'aax
$ws.Range("H136").Value = 2344.25
$ws.Range("I136").Value = 1789.9584
$ws.Range("J136").Value = 4007.125
$ws.Range("K136").Value = 5369.8752
$ws.Range("L136").Value = 12021.375
$ws.Range("M136").Value = -2819.8752
$ws.Range("N136").Value = -17121.375

$ws = $wb.Worksheets.Item("WVR")
# WVR row 16: Keep It under Wraps
$ws.Range("H16").Value = 46803.668
$ws.Range("J16").Value = 46803.668
$ws.Range("L16").Value = 46803.668
$ws.Range("N16").Value = -47387.668

# WVR row 119: A Job Well Done
$ws.Range("H119").Value = 45396.332
$ws.Range("J119").Value = 45396.332
$ws.Range("L119").Value = 45396.332
$ws.Range("N119").Value = -55072.332

# WVR row 120: A Turban for the Ages
$ws.Range("H120").Value = 42010.668
$ws.Range("J120").Value = 42010.668
$ws.Range("L120").Value = 42010.668
$ws.Range("N120").Value = -51686.668

# WVR row 137: Traditional Trousers
$ws.Range("H137").Value = 22384
$ws.Range("J137").Value = 22384
$ws.Range("L137").Value = 22384
$ws.Range("N137").Value = -32584
